# Weekly update: insert a new price record (row 58) for
# "Vega Modelo de Temuco" / Achicoria, pushing the existing rows 58-61
# down to 59-62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 58 (existing rows 58..61 shift to 59..62,
# inheriting their formatting/style automatically).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the latest week's data.
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 44783
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = 100112010
$ws.Range("G58").Value = "Achicoria"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 100
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 10000
$ws.Range("M58").Value = 10000
$ws.Range("N58").Value = "`$/caja 18 unidades"
$ws.Range("O58").Value = "Región Metropolitana"
$ws.Range("P58").Value = 556
$ws.Range("Q58").Value = 18
$ws.Range("R58").Value = "Hortaliza"
